# Update need_to_buy.xlsx values in columns C (fcs) and F (need_to_buy)
# per refreshed R computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 11712.2185924037
$ws.Range("F2").Value = -17.2722172767997

$ws.Range("C3").Value = 10820.7120890551
$ws.Range("F3").Value = 292.438660746478

$ws.Range("C4").Value = 7327.1569982774
$ws.Range("F4").Value = 125.213462390191

$ws.Range("C5").Value = 7267.36590227021
$ws.Range("F5").Value = 121.277490234164

$ws.Range("C6").Value = 11254.0169659754
$ws.Range("F6").Value = 303.220202389622

$ws.Range("C7").Value = 12211.7102819467
$ws.Range("F7").Value = 366.349600847913

$ws.Range("C9").Value = 12375.2161535474
$ws.Range("F9").Value = 373.162345497944

$ws.Range("C10").Value = 11847.992695046
$ws.Range("F10").Value = 351.194701393715

$ws.Range("C11").Value = 8576.40814513318
$ws.Range("F11").Value = 198.497596841317

$ws.Range("C12").Value = 8733.33190970265
$ws.Range("F12").Value = 204.713397423308

$ws.Range("C13").Value = 12742.5247986484
$ws.Range("F13").Value = 374.16236451889

$ws.Range("C14").Value = 12941.6823862053
$ws.Range("F14").Value = 382.460597333763

$ws.Range("C15").Value = 12776.5363522979
$ws.Range("F15").Value = 375.579512587619
